# Apply price/volume updates to the cryptos worksheet.
# Values are entered with a leading apostrophe to force Excel to store them
# as literal text (matching the inline-string cells in the source workbook)
# instead of auto-converting numeric-looking text into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.150.37"
$ws.Range("E2").Value = "'  -0.47%  "
$ws.Range("D3").Value = "'2.526.80"
$ws.Range("E3").Value = "'  +0.57%  "
$ws.Range("D5").Value = "'536.58"
$ws.Range("E5").Value = "'  -0.92%  "
$ws.Range("D6").Value = "'137.13"
$ws.Range("E6").Value = "'  -1.74%  "
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("E8").Value = "'  +0.80%  "
$ws.Range("D9").Value = "'2.527.22"
$ws.Range("E9").Value = "'  +0.50%  "
$ws.Range("E10").Value = "'  -0.29%  "
$ws.Range("E11").Value = "'  -1.94%  "
$ws.Range("E12").Value = "'  -1.12%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "'  -1.56%  "
$ws.Range("D14").Value = "'2.974.78"
$ws.Range("E14").Value = "'  +0.38%  "
$ws.Range("D15").Value = "'23.06"
$ws.Range("E15").Value = "'  -1.37%  "
$ws.Range("D16").Value = "'59.036.58"
$ws.Range("E16").Value = "'  -0.47%  "
$ws.Range("E17").Value = "'  -1.36%  "
$ws.Range("D18").Value = "'2.528.95"
$ws.Range("E18").Value = "'  +0.74%  "
$ws.Range("D19").Value = "'11.15"
$ws.Range("E19").Value = "'  +0.61%  "
$ws.Range("D20").Value = "'4.28"
$ws.Range("E20").Value = "'  -0.33%  "
$ws.Range("D21").Value = "'323.53"
$ws.Range("E21").Value = "'  -0.65%  "
$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "'  +1.31%  "
$ws.Range("D24").Value = "'65.83"
$ws.Range("E24").Value = "'  +3.70%  "
$ws.Range("D25").Value = "'0.423"
$ws.Range("E25").Value = "'  +0.06%  "
$ws.Range("E26").Value = "'  -1.67%  "
$ws.Range("E27").Value = "'  +0.03%  "
$ws.Range("D28").Value = "'7.54"
$ws.Range("E28").Value = "'  -3.39%  "
$ws.Range("D29").Value = "'0.0₃0775"
$ws.Range("E29").Value = "'  -0.59%  "
$ws.Range("D30").Value = "'6.71"
$ws.Range("E30").Value = "'  -1.89%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "'  -1.38%  "
$ws.Range("D32").Value = "'167.51"
$ws.Range("E32").Value = "'  +2.48%  "
$ws.Range("E33").Value = "'  +5.10%  "
$ws.Range("E34").Value = "'  +0.02%  "
$ws.Range("E35").Value = "'  +1.65%  "
$ws.Range("E36").Value = "'  -0.29%  "
$ws.Range("D37").Value = "'4.11"
$ws.Range("E37").Value = "'  -2.75%  "
$ws.Range("E38").Value = "'  -3.11%  "
$ws.Range("D39").Value = "'36.72"
$ws.Range("E39").Value = "'  -0.59%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = "'  +0.32%  "
$ws.Range("D41").Value = "'3.62"
$ws.Range("E41").Value = "'  -1.66%  "
$ws.Range("D42").Value = "'285.64"
$ws.Range("E42").Value = "'  +1.46%  "
$ws.Range("E43").Value = "'  -1.75%  "
$ws.Range("D44").Value = "'132.38"
$ws.Range("E44").Value = "'  +5.63%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.18%  "
$ws.Range("D46").Value = "'0.607"
$ws.Range("E46").Value = "'  +1.64%  "
$ws.Range("E47").Value = "'  +0.33%  "
$ws.Range("D48").Value = "'0.0925"
$ws.Range("E48").Value = "'  -1.16%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "'  -1.00%  "
$ws.Range("E50").Value = "'  -1.58%  "
$ws.Range("D51").Value = "'17.38"
$ws.Range("E51").Value = "'  -2.78%  "
